$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.990.21"
$ws.Range("E2").Value = "  +2.99%  "

$ws.Range("D3").Value = "3.237.07"
$ws.Range("E3").Value = "  +7.03%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'579.94"
$ws.Range("E5").Value = "  +5.17%  "

$ws.Range("D6").Value = "'151.26"
$ws.Range("E6").Value = "  +8.92%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.227.06"
$ws.Range("E8").Value = "  +6.92%  "

$ws.Range("D9").Value = "'0.512"
$ws.Range("E9").Value = "  +6.27%  "

$ws.Range("D10").Value = "'7.06"
$ws.Range("E10").Value = "  +10.88%  "

$ws.Range("D11").Value = "'0.162"
$ws.Range("E11").Value = "  +6.98%  "

$ws.Range("D12").Value = "'0.486"
$ws.Range("E12").Value = "  +6.35%  "

$ws.Range("D13").Value = "'37.56"
$ws.Range("E13").Value = "  +4.32%  "

$ws.Range("D14").Value = "'0.0000233"
$ws.Range("E14").Value = "  +7.30%  "

$ws.Range("D15").Value = "3.751.73"
$ws.Range("E15").Value = "  +7.14%  "

$ws.Range("D16").Value = "66.074.00"
$ws.Range("E16").Value = "  +3.11%  "

$ws.Range("D17").Value = "'543.71"
$ws.Range("E17").Value = "  +13.67%  "

$ws.Range("D18").Value = "3.240.71"
$ws.Range("E18").Value = "  +7.25%  "

$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("D20").Value = "'7.09"
$ws.Range("E20").Value = "  +7.37%  "

$ws.Range("D21").Value = "'14.47"
$ws.Range("E21").Value = "  +7.59%  "

$ws.Range("D22").Value = "'0.742"
$ws.Range("E22").Value = "  +9.55%  "

$ws.Range("D23").Value = "'7.85"
$ws.Range("E23").Value = "  +11.78%  "

$ws.Range("D24").Value = "'13.43"
$ws.Range("E24").Value = "  +7.32%  "

$ws.Range("D25").Value = "'80.80"
$ws.Range("E25").Value = "  +3.87%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +19.82%  "

$ws.Range("E28").Value = "  +10.39%  "

$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  +7.30%  "

$ws.Range("D30").Value = "'27.59"
$ws.Range("E30").Value = "  +7.59%  "

$ws.Range("D31").Value = "'2.74"
$ws.Range("E31").Value = "  +5.95%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("E33").Value = "  +6.97%  "

$ws.Range("D34").Value = "'568.93"
$ws.Range("E34").Value = "  +9.84%  "

$ws.Range("D35").Value = "'5.63"
$ws.Range("E35").Value = "  +5.02%  "

$ws.Range("D36").Value = "'6.32"
$ws.Range("E36").Value = "  +7.66%  "

$ws.Range("D37").Value = "'55.14"
$ws.Range("E37").Value = "  +5.85%  "

$ws.Range("D38").Value = "'0.0453"
$ws.Range("E38").Value = "  +14.32%  "

$ws.Range("D39").Value = "'0.0858"
$ws.Range("E39").Value = "  +8.65%  "

$ws.Range("E40").Value = "  +6.51%  "

$ws.Range("D41").Value = "3.193.24"
$ws.Range("E41").Value = "  +11.26%  "

$ws.Range("D42").Value = "'2.94"
$ws.Range("E42").Value = "  +9.21%  "

$ws.Range("D43").Value = "'8.56"
$ws.Range("E43").Value = "  +4.50%  "

$ws.Range("D44").Value = "'0.281"
$ws.Range("E44").Value = "  +17.51%  "

$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  +11.78%  "

$ws.Range("D46").Value = "'26.35"
$ws.Range("E46").Value = "  +7.14%  "

$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("E48").Value = "  +6.86%  "

$ws.Range("D49").Value = "'126.01"
$ws.Range("E49").Value = "  +5.66%  "

$ws.Range("E50").Value = "  +4.37%  "

$ws.Range("E51").Value = "  +9.55%  "
